# Updates the cryptos price/volume list to the latest scraped values.
# Rows 35/36 (Mantle <-> PEPE) and 39/40 (OKB <-> Cosmos) swapped their
# ranking position, so their Coin/Link/Price/Volume columns are rewritten
# in place rather than moved.
#
# Column D (Price) is forced to text format before assignment so that
# Excel does not reinterpret dotted/thousands values (e.g. "62.209.97")
# or trailing-zero decimals (e.g. "10.00") as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.209.97"
$ws.Range("E2").Value = "  -2.36%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.000.20"
$ws.Range("E3").Value = "  -2.35%  "

$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.66"
$ws.Range("E5").Value = "  -1.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.79"
$ws.Range("E6").Value = "  -5.77%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.525"
$ws.Range("E8").Value = "  -2.51%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.998.05"
$ws.Range("E9").Value = "  -2.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.148"
$ws.Range("E10").Value = "  -5.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.81"
$ws.Range("E11").Value = "  -1.92%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("E12").Value = "  +0.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000228"
$ws.Range("E13").Value = "  -3.98%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.40"
$ws.Range("E14").Value = "  -6.56%  "

$ws.Range("E15").Value = "  +1.88%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.490.40"
$ws.Range("E16").Value = "  -2.50%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.09"
$ws.Range("E17").Value = "  -1.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.128.69"
$ws.Range("E18").Value = "  -2.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.003.17"
$ws.Range("E19").Value = "  -2.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "460.39"
$ws.Range("E20").Value = "  -4.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.95"
$ws.Range("E21").Value = "  -4.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.685"
$ws.Range("E22").Value = "  -3.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.43"
$ws.Range("E23").Value = "  -2.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.34"
$ws.Range("E24").Value = "  -0.93%  "

$ws.Range("E25").Value = "  -10.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.26"
$ws.Range("E26").Value = "  -5.04%  "

$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.00"
$ws.Range("E28").Value = "  -6.90%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("E30").Value = "  -3.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.01"
$ws.Range("E31").Value = "  -6.34%  "

$ws.Range("E32").Value = "  -7.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.10"
$ws.Range("E33").Value = "  +2.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.108"
$ws.Range("E34").Value = "  -3.73%  "

$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0794"
$ws.Range("E35").Value = "  -3.18%  "

$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.02"
$ws.Range("E36").Value = "  -3.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.74"
$ws.Range("E37").Value = "  -5.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.10"
$ws.Range("E38").Value = "  -6.03%  "

$ws.Range("B39").Value = "Cosmos"
$ws.Range("C39").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.18"
$ws.Range("E39").Value = "  -0.87%  "

$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.27"
$ws.Range("E40").Value = "  -0.77%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.87"
$ws.Range("E41").Value = "  -11.28%  "

$ws.Range("E42").Value = "  +1.35%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "393.03"
$ws.Range("E43").Value = "  -10.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0357"
$ws.Range("E44").Value = "  -2.41%  "

$ws.Range("E45").Value = "  -7.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.720.85"
$ws.Range("E46").Value = "  -4.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.65"
$ws.Range("E47").Value = "  -6.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.17"
$ws.Range("E48").Value = "  -2.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.109"
$ws.Range("E50").Value = "  -0.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.19"
$ws.Range("E51").Value = "  -2.79%  "
